$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "320018813081"
$ws.Range("P2").ClearFormats()
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "320018813092"
$ws.Range("P3").ClearFormats()
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "320018813129"
$ws.Range("P4").ClearFormats()
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "320018813140"
$ws.Range("P5").ClearFormats()
$ws.Range("P6").NumberFormat = "@"
$ws.Range("P6").Value = "320018813184"
$ws.Range("P6").ClearFormats()
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "320018813200"
$ws.Range("P7").ClearFormats()
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "320018813232"
$ws.Range("P8").ClearFormats()
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "320018813254"
$ws.Range("P9").ClearFormats()
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "320018813287"
$ws.Range("P10").ClearFormats()
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "320018813302"
$ws.Range("P11").ClearFormats()
$ws.Range("P12").NumberFormat = "@"
$ws.Range("P12").Value = "320018813346"
$ws.Range("P12").ClearFormats()
$ws.Range("P13").NumberFormat = "@"
$ws.Range("P13").Value = "320018813449"
$ws.Range("P13").ClearFormats()
$ws.Range("P14").NumberFormat = "@"
$ws.Range("P14").Value = "320018813471"
$ws.Range("P14").ClearFormats()
$ws.Range("P15").NumberFormat = "@"
$ws.Range("P15").Value = "320018813493"
$ws.Range("P15").ClearFormats()
$ws.Range("P16").NumberFormat = "@"
$ws.Range("P16").Value = "320018813520"
$ws.Range("P16").ClearFormats()
$ws.Range("P17").NumberFormat = "@"
$ws.Range("P17").Value = "320018813541"
$ws.Range("P17").ClearFormats()
$ws.Range("P18").NumberFormat = "@"
$ws.Range("P18").Value = "320018813585"
$ws.Range("P18").ClearFormats()
$ws.Range("P19").NumberFormat = "@"
$ws.Range("P19").Value = "320018813600"
$ws.Range("P19").ClearFormats()
$ws.Range("P20").NumberFormat = "@"
$ws.Range("P20").Value = "320018813850"
$ws.Range("P20").ClearFormats()
$ws.Range("P21").NumberFormat = "@"
$ws.Range("P21").Value = "320018813872"
$ws.Range("P21").ClearFormats()
$ws.Range("P22").NumberFormat = "@"
$ws.Range("P22").Value = "320018813909"
$ws.Range("P22").ClearFormats()
$ws.Range("P23").NumberFormat = "@"
$ws.Range("P23").Value = "320018813910"
$ws.Range("P23").ClearFormats()
$ws.Range("P24").NumberFormat = "@"
$ws.Range("P24").Value = "320018813920"
$ws.Range("P24").ClearFormats()
$ws.Range("P25").NumberFormat = "@"
$ws.Range("P25").Value = "320018813931"
$ws.Range("P25").ClearFormats()
$ws.Range("P26").NumberFormat = "@"
$ws.Range("P26").Value = "320018813942"
$ws.Range("P26").ClearFormats()
